# Daily attendance processing - 2026-01-20 11:39:23
# Reorders the comma-separated "Recorded By" values in column G for the
# affected rows (reversing the order of the listed names/emails).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, System, system"
    3   = "System, dnasr281@gmail.com"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    6   = "System, dnasr281@gmail.com"
    8   = "backup@backdoor.com, System"
    28  = "backup@backdoor.com, System, system"
    29  = "System, dnasr281@gmail.com"
    30  = "backup@backdoor.com, System"
    31  = "backup@backdoor.com, System"
    32  = "System, dnasr281@gmail.com"
    34  = "backup@backdoor.com, System"
    54  = "backup@backdoor.com, System, system"
    55  = "System, dnasr281@gmail.com"
    56  = "backup@backdoor.com, System"
    57  = "backup@backdoor.com, System"
    58  = "System, dnasr281@gmail.com"
    60  = "backup@backdoor.com, System"
    80  = "backup@backdoor.com, System"
    81  = "backup@backdoor.com, System"
    82  = "backup@backdoor.com, System"
    87  = "admin@admin.com, dnasr281@gmail.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    113 = "admin@admin.com, dnasr281@gmail.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    139 = "admin@admin.com, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
